# Add 2022-Q4 data:
#  - Insert a new "2022-Q4" worksheet (cloned from "2022-Q3" so it keeps the
#    same layout/styles) right after "总计" and before "2022-Q3".
#  - Fill it with the new quarter's fund-holding figures.
#  - Update the "总计" (totals) sheet with a new top row for 2022-Q4 and
#    shift the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the "2022-Q4" sheet by copying "2022-Q3" (same columns/styles),
#    placed immediately before "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Row 2 - fund 005947
$q4.Range("D2").Value = "'0.50"
$q4.Range("E2").Value = "'90.73"
$q4.Range("F2").Value = "'7.80"
$q4.Range("G2").Value = "'0.0390"
$q4.Range("H2").Value = 2

# Row 3 - fund 005948
$q4.Range("D3").Value = "'0.09"
$q4.Range("E3").Value = "'90.73"
$q4.Range("F3").Value = "'7.80"
$q4.Range("G3").Value = "'0.0070"
$q4.Range("H3").Value = 2

# ---------------------------------------------------------------------
# 2) Update the "总计" totals sheet: add a 2022-Q4 row at the top of the
#    data (row 2) and push Q3/Q2/Q1 down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give the new bottom row (row 5) the same index-column style as row 4.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

# 2022-Q1 moves from row 4 to row 5
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.05

# 2022-Q2 moves from row 3 to row 4
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.04

# 2022-Q3 moves from row 2 to row 3
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.04

# 2022-Q4 becomes the new row 2
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.05

# ---------------------------------------------------------------------
# 3) Restore the original active sheet ("2022-Q1" was the selected tab
#    before this edit; inserting/copying sheets shifts the active tab to
#    the freshly-created one, so put the selection back where it was).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
